# Apply updated crypto price/volume data to Sheet1.
# Column D holds price strings that can look numeric (e.g. "606.07" or
# thousand-dotted "63.166.10"); Excel's COM Range.Value setter auto-detects
# plain-numeric-looking text and silently coerces it to a Number (dropping
# formatting like trailing zeros). Forcing NumberFormat="@" (Text) on those
# cells first keeps them as text, matching the original inlineStr cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.166.10"
$ws.Range("E2").Value = "  -4.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.101.81"
$ws.Range("E3").Value = "  -4.80%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.07"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.29"
$ws.Range("E6").Value = "  -8.96%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.095.78"
$ws.Range("E8").Value = "  -4.84%  "
$ws.Range("E9").Value = "  -4.50%  "
$ws.Range("E10").Value = "  -7.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.20"
$ws.Range("E11").Value = "  -8.87%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.465"
$ws.Range("E12").Value = "  -5.91%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000247"
$ws.Range("E13").Value = "  -8.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.85"
$ws.Range("E14").Value = "  -9.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.619.20"
$ws.Range("E15").Value = "  -4.64%  "
$ws.Range("E16").Value = "  +1.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.306.32"
$ws.Range("E17").Value = "  -4.28%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.109.56"
$ws.Range("E18").Value = "  -4.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.73"
$ws.Range("E19").Value = "  -7.93%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "471.61"
$ws.Range("E20").Value = "  -5.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.45"
$ws.Range("E21").Value = "  -5.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.702"
$ws.Range("E22").Value = "  -6.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.64"
$ws.Range("E23").Value = "  -5.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.34"
$ws.Range("E24").Value = "  -8.82%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.30"
$ws.Range("E25").Value = "  -3.53%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.75"
$ws.Range("E27").Value = "  -9.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.28"
$ws.Range("E28").Value = "  -9.86%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.07"
$ws.Range("E29").Value = "  -12.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.72"
$ws.Range("E30").Value = "  -4.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.112"
$ws.Range("E31").Value = "  -15.75%  "
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.64"
$ws.Range("E33").Value = "  -7.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.98"
$ws.Range("E34").Value = "  -6.57%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.10"
$ws.Range("E35").Value = "  -4.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.88"
$ws.Range("E36").Value = "  -8.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.59"
$ws.Range("E37").Value = "  -5.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0731"
$ws.Range("E38").Value = "  -7.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "448.72"
$ws.Range("E39").Value = "  -10.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.89"
$ws.Range("E40").Value = "  -16.85%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0387"
$ws.Range("E41").Value = "  -8.38%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.116"
$ws.Range("E42").Value = "  -9.53%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.24"
$ws.Range("E43").Value = "  -5.99%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.825.95"
$ws.Range("E44").Value = "  -5.56%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.261"
$ws.Range("E45").Value = "  -9.94%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.23"
$ws.Range("E46").Value = "  -13.80%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.38"
$ws.Range("E47").Value = "  -3.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.998"
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.67"
$ws.Range("E49").Value = "  -10.77%  "
$ws.Range("E50").Value = "  -5.79%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "117.82"
$ws.Range("E51").Value = "  -2.66%  "
